# Applies the text edits described by the commit diff to the active presentation.
#
# Technique notes (this COM host's quirks, discovered by experimentation):
#  - TextRange.Paragraphs(n,1) returns the whole paragraph (it is a single run in
#    every paragraph we touch here). Assigning .Text directly on a *freshly
#    fetched* paragraph range keeps rPr and collapses to one run, matching the
#    target XML, but only if the new text shares no run-splitting "common
#    prefix" with what is already there — so we first stomp the paragraph with
#    an unrelated placeholder, then set the real text (also freshly fetched).
#  - Caching a Paragraphs()/Characters() range object across more than one
#    .Text assignment silently no-ops the second write, so every assignment
#    below re-fetches the range right before it is used.
#  - TextRange.Characters(start, length) can target an exact existing run's
#    span (by 1-based character offset) and rewrite just that run in place
#    without disturbing neighboring (e.g. bold) runs.

$p = $ppt.ActivePresentation

# --- Slide 4 ("Data-  Town Centers") : Content Placeholder 2, paragraph 1 ---
$s4 = $p.Slides.Item(4)
$tr4 = $s4.Shapes.Item(2).TextFrame.TextRange
$tr4.Paragraphs(1,1).Text = "X"
$tr4.Paragraphs(1,1).Text = "We identified each town center’s GPS location and then evaluated the restaurants at 1000m, 2000m and 3000m from the known town centers.   2000m should capture the restaurants within walking distance of the town centers.  Googles Geocoding API will be used for identifying the established latitude and longitude of each town center.  "

# --- Slide 5 ("Methodology") : Content Placeholder 2 ---
$s5 = $p.Slides.Item(5)
$tr5 = $s5.Shapes.Item(2).TextFrame.TextRange
# paragraph 2: drop trailing clause about the combined suburb counts
$tr5.Paragraphs(2,1).Text = "X"
$tr5.Paragraphs(2,1).Text = "We then evaluated the type and number of food venues to determine how many of each were in each suburb."
# paragraph 6: rewrite the word-cloud sentence ending
$tr5.Paragraphs(6,1).Text = "X"
$tr5.Paragraphs(6,1).Text = "In order to illustrate the results, we created 3 word clouds to evaluate varying distances from the town center."

# --- Slide 7 ("Analysis at 1000 meters") : Content Placeholder 2, first run ---
$s7 = $p.Slides.Item(7)
$tr7 = $s7.Shapes.Item(2).TextFrame.TextRange
$old7 = "Our initial study at 1000 meters from town centers, we found that the most common venues were "
$idx7 = $tr7.Text.IndexOf($old7)
$tr7.Characters($idx7 + 1, $old7.Length).Text = "In our initial study at 1000 meters from town centers, we found that the most common venues were "

# --- Slide 8 ("Analysis at 2000 meters") : Content Placeholder 2, last run ---
$s8 = $p.Slides.Item(8)
$tr8 = $s8.Shapes.Item(2).TextFrame.TextRange
$old8 = " restaurant the farther you go from town centers."
$idx8 = $tr8.Text.IndexOf($old8)
$tr8.Characters($idx8 + 1, $old8.Length).Text = " restaurants the farther you go from town centers."

# --- Slide 10 ("Conclusion") : Content Placeholder 2, paragraph 3 (drop stray leading space) ---
$s10 = $p.Slides.Item(10)
$tr10 = $s10.Shapes.Item(2).TextFrame.TextRange
$tr10.Paragraphs(3,1).Text = "X"
$tr10.Paragraphs(3,1).Text = "We can also conclude that American, Mexican and Pizza can be a success regardless of distance to the town center.  "

# --- Presentation-level: touch the (empty) slide-guide extension list, mirroring the
#     no-op bookkeeping PowerPoint writes to presentation.xml on save. Best effort —
#     harmless if the host does not model Guides persistence.
try {
    $guides = $p.Guides
    $guides.Add(1, 3000) | Out-Null
} catch {
}
